# fix: Moved env variables in pipeline
#
# Adds three reference-link paragraphs right after the "References:" heading
# (and before the existing blank paragraph that follows it), each paragraph
# containing a single hyperlink styled with the built-in "Hyperlink" style.

$d = $word.ActiveDocument

# Locate the "References:" heading paragraph and get a collapsed range right
# at its end (= right before the blank paragraph that already follows it).
$refRange = $d.Content
$refRange.Find.Execute("References:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$refPara = $refRange.Paragraphs(1)
$insertPoint = $d.Range($refPara.Range.End, $refPara.Range.End)

# Use unique sentinel markers as plain text first, separated by paragraph
# marks ("`r"). This creates three clean new paragraphs (no stray empty
# runs) that we then convert to hyperlinks one at a time below.
$marker1 = "___REF_LINK_1___"
$marker2 = "___REF_LINK_2___"
$marker3 = "___REF_LINK_3___"
$insertPoint.InsertAfter($marker1 + "`r" + $marker2 + "`r" + $marker3)

$links = @(
    @{ Marker = $marker1; Url = "https://docs.github.com/en/actions/automating-builds-and-tests/building-and-testing-net"; Text = "Building and testing .NET - GitHub Docs" },
    @{ Marker = $marker2; Url = "https://learn.microsoft.com/en-us/aspnet/core/tutorials/min-web-api"; Text = "Tutorial: Create a minimal web API with ASP.NET Core | Microsoft Learn" },
    @{ Marker = $marker3; Url = "https://docs.github.com/en/actions/learn-github-actions/variables"; Text = "Environment variables - GitHub Docs" }
)

foreach ($link in $links) {
    $found = $d.Content
    $found.Find.Execute($link.Marker, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $d.Hyperlinks.Add($found, $link.Url, "", "", $link.Text)
}

"Inserted " + $links.Count + " reference hyperlinks after 'References:'"
